# feat: add 2022-Q4 data
#
# Original workbook has sheets: 总计, 2022-Q3, 2022-Q2
# Target workbook has sheets:  总计, 2022-Q4, 2022-Q3, 2022-Q2
#  - A new "2022-Q4" sheet is inserted right after "总计", with fresh data.
#  - The old "2022-Q3" / "2022-Q2" sheets keep their data unchanged, just shift position.
#  - The "总计" (summary) sheet gains a row for 2022-Q4 and keeps the Q3/Q2 rows (shifted).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right before the existing "2022-Q3"
#    sheet (i.e. as the 2nd sheet, right after "总计").
#
# NOTE: worksheet object references in this engine track *position*, not
# identity, so once Worksheets.Add() shifts sheets around, any reference
# obtained beforehand now points at whatever sheet occupies that old slot.
# To stay safe we re-fetch sheets *by name* right after the insertion, and
# again at the very end (right before Activate()).
# ---------------------------------------------------------------------------
$q3SheetBeforeInsert = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($q3SheetBeforeInsert)
$q4Sheet.Name = "2022-Q4"

# Re-fetch stable references now that the sheet collection is final.
$sumSheet = $wb.Worksheets.Item("总计")
$q3Sheet  = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as genuine text (even if it "looks" like
# a number), without leaving a stray NumberFormat-driven cell style behind.
# Trick: flip the cell to text format, assign the value, then restore the
# cell's original formatting (copied from a plain/no-special-style cell) so
# the stored value stays textual but the visual style/index goes back to
# whatever it was (no extra style entries created).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($ws, $addr, $val, $plainFormatCell)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $plainFormatCell.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 2) Build up the "2022-Q4" sheet: copy the header row + per-row formatting
#    template from "2022-Q3" (same column layout/styles), then fill in the
#    2022-Q4 values (17 funds).
# ---------------------------------------------------------------------------

# Header row (values + formats are identical across quarters)
$q3Sheet.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$q3Sheet.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4163)

# Row formatting template for the 13 rows 2022-Q3 already has
$q3Sheet.Range("A2:H14").Copy()
$q4Sheet.Range("A2:H14").PasteSpecial(-4122)

# 2022-Q4 needs 17 data rows (rows 2-18); clone the formatting of the last
# template row (row 14) down to the extra rows 15-18.
$q3Sheet.Range("A14:H14").Copy()
$q4Sheet.Range("A15:H15").PasteSpecial(-4122)
$q4Sheet.Range("A16:H16").PasteSpecial(-4122)
$q4Sheet.Range("A17:H17").PasteSpecial(-4122)
$q4Sheet.Range("A18:H18").PasteSpecial(-4122)

# A "plain" (no special style) cell to re-copy formatting from whenever we
# need to strip the stray style that NumberFormat="@" leaves behind.
$plainCell = $q4Sheet.Range("B2")

$q4Data = @(
    @("159941", "广发纳斯达克100ETF（QDII）", "114.77", "90.42", "2.21", "2.5364", 8),
    @("513100", "国泰纳斯达克100（QDII-ETF）", "51.50", "90.72", "2.33", "1.2000", 7),
    @("159632", "华安纳斯达克100ETF（QDII）", "30.39", "87.77", "2.18", "0.6625", 8),
    @("160213", "国泰纳斯达克100指数（QDII）", "15.65", "90.80", "2.25", "0.3521", 8),
    @("000834", "大成纳斯达克100指数（QDII）", "15.61", "81.77", "2.06", "0.3216", 8),
    @("161130", "易方达纳斯达克100指数人民币（QDII-LOF）", "7.77", "90.34", "2.24", "0.1740", 8),
    @("003722", "易方达纳斯达克100指数美元（QDII-LOF）A", "7.77", "90.34", "2.24", "0.1740", 8),
    @("016532", "嘉实纳斯达克100指数（QDII）A人民币", "1.12", "94.67", "2.35", "0.0263", 8),
    @("016533", "嘉实纳斯达克100指数（QDII）C人民币", "1.12", "94.67", "2.35", "0.0263", 8),
    @("016534", "嘉实纳斯达克100指数（QDII）A美元现汇", "1.12", "94.67", "2.35", "0.0263", 8),
    @("016535", "嘉实纳斯达克100指数（QDII）C美元现汇", "1.12", "94.67", "2.35", "0.0263", 8),
    @("016055", "博时纳斯达克100指数（QDII）A人民币", "1.06", "90.62", "2.25", "0.0238", 8),
    @("016057", "博时纳斯达克100指数（QDII）C人民币", "1.06", "90.62", "2.25", "0.0238", 8),
    @("016056", "博时纳斯达克100指数（QDII）A美元现汇", "1.06", "90.62", "2.25", "0.0238", 8),
    @("016058", "博时纳斯达克100指数（QDII）C美元现汇", "1.06", "90.62", "2.25", "0.0238", 8),
    @("012870", "易方达纳斯达克100指数人民币（QDII-LOF）C", "0.21", "90.34", "2.24", "0.0047", 8),
    @("012871", "易方达纳斯达克100指数美元（QDII-LOF）C", "0.21", "90.34", "2.24", "0.0047", 8)
)

$row = 2
foreach ($rec in $q4Data) {
    $q4Sheet.Range("A$row").Value = ($row - 2)
    Set-TextValue $q4Sheet "B$row" $rec[0] $plainCell
    Set-TextValue $q4Sheet "C$row" $rec[1] $plainCell
    Set-TextValue $q4Sheet "D$row" $rec[2] $plainCell
    Set-TextValue $q4Sheet "E$row" $rec[3] $plainCell
    Set-TextValue $q4Sheet "F$row" $rec[4] $plainCell
    Set-TextValue $q4Sheet "G$row" $rec[5] $plainCell
    $q4Sheet.Range("H$row").Value = $rec[6]
    $row = $row + 1
}

# Match page margins used by the other worksheets (0.75/0.75/1/1/0.5/0.5 in).
$q4Sheet.PageSetup.LeftMargin = 54
$q4Sheet.PageSetup.RightMargin = 54
$q4Sheet.PageSetup.TopMargin = 72
$q4Sheet.PageSetup.BottomMargin = 72
$q4Sheet.PageSetup.HeaderMargin = 36
$q4Sheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: row 2 becomes 2022-Q4, row 3 becomes
#    2022-Q3 (same numbers the old row 2 had), and a new row 4 is added for
#    2022-Q2 (same numbers the old row 3 had).
# ---------------------------------------------------------------------------
$sumSheet.Range("A3:D3").Copy()
$sumSheet.Range("A4:D4").PasteSpecial(-4122)

$sumSheet.Range("B2").Value = "2022-Q4"
$sumSheet.Range("C2").Value = 17
$sumSheet.Range("D2").Value = 5.63

$sumSheet.Range("B3").Value = "2022-Q3"
$sumSheet.Range("C3").Value = 13
$sumSheet.Range("D3").Value = 6.91

$sumSheet.Range("A4").Value = 2
$sumSheet.Range("B4").Value = "2022-Q2"
$sumSheet.Range("C4").Value = 13
$sumSheet.Range("D4").Value = 7.26

# ---------------------------------------------------------------------------
# 4) Restore the original "active sheet" (2022-Q2 was tabSelected before the
#    edit); adding the new sheet made it active, so reselect 2022-Q2.
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Activate()
